# semana 41 de 2025
# Adds week-41 column (AR) data to the weekly IRA-Ext report sheet,
# plus a couple of late-arriving corrections to existing weeks
# (AQ44, and several values in row 36 for columns G..U).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header for week 41 (stored as text, like the other week headers) ---
$ws.Range("AR1").Value = "'41"

# --- New week-41 values (column AR) per UPGD row ---
$ar = @{
    2  = 59
    5  = 0
    6  = 88
    7  = 23
    8  = 18
    9  = 5
    11 = 1
    13 = 2
    14 = 1
    16 = 1
    17 = 2
    22 = 2
    23 = 7
    24 = 1
    25 = 23
    26 = 2
    29 = 0
    30 = 60
    31 = 3
    35 = 68
    36 = 2
    37 = 3
    38 = 83
    41 = 5
    42 = 12
    43 = 19
    44 = 133
    45 = 68
    46 = 88
    47 = 0
    48 = 104
    49 = 5
    50 = 0
    51 = 11
    53 = 11
    54 = 0
    55 = 0
    56 = 10
    57 = 39
    58 = 16
}

foreach ($row in $ar.Keys) {
    $ws.Range("AR$row").Value = $ar[$row]
}

# --- Late correction: row 44 also gets a week-40 (AQ) value now ---
$ws.Range("AQ44").Value = 165

# --- Corrections to row 36, weeks 4-7 and 14-18 (columns G,H,I,J,L,Q,R,S,T,U) ---
$ws.Range("G36").Value = 2
$ws.Range("H36").Value = 3
$ws.Range("I36").Value = 3
$ws.Range("J36").Value = 5
$ws.Range("L36").Value = 3
$ws.Range("Q36").Value = 6
$ws.Range("R36").Value = 1
$ws.Range("S36").Value = 3
$ws.Range("T36").Value = 2
$ws.Range("U36").Value = 3
